$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 0.1 -> 1.0 (leading apostrophe forces text so "1.0" isn't
# auto-converted to the number 1)
$ws.Range("D2").Value = "'1.0"

# TC1 (rows 9-13): step 1 text capitalization fix
$ws.Range("B10").Value = "Usuario do Sistema inicia a tela de login atraves da opcao de Login no canto superior direito"
# TC1: swap step 2 and step 3 texts
$ws.Range("B11").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$ws.Range("D11").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"
$ws.Range("B12").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"

# TC2 (rows 19-23): step 1 text capitalization fix
$ws.Range("B20").Value = "Usuario do Sistema inicia a tela de login atraves da opcao de Login no canto superior direito"
# TC2: step 2 text updated
$ws.Range("B21").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"

# TC3 (rows 29-33): step 1 text capitalization fix
$ws.Range("B30").Value = "Usuario do Sistema inicia a tela de login atraves da opcao de Login no canto superior direito"
# TC3: step 2 expected result updated
$ws.Range("D31").Value = "SYSTEM alerta que o CAS (sistema de autorizacao login-senha) esta fora do ar"
